$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Move "Line_DLR" sheet so it becomes the 2nd sheet (right after
#    "Stations" and before "Line_Central Line").
# ------------------------------------------------------------------
$stations = $wb.Worksheets.Item("Stations")
$dlr = $wb.Worksheets.Item("Line_DLR")
$dlr.Move($null, $stations)

# ------------------------------------------------------------------
# 2. Fill in missing "opened" year (column B) for DLR stations on the
#    "Stations" sheet that previously had a blank year.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Stations")
$ws.Range("B150").Value = 1987   # Devons Road
$ws.Range("B298").Value = 2005   # King George V
$ws.Range("B310").Value = 2007   # Langdon Park
$ws.Range("B319").Value = 1840   # Limehouse
$ws.Range("B322").Value = 2005   # London City Airport
$ws.Range("B402").Value = 2005   # Pontoon Dock
$ws.Range("B403").Value = 1987   # Poplar
$ws.Range("B443").Value = 1987   # Shadwell
$ws.Range("B486").Value = 2011   # Star Lane
$ws.Range("B493").Value = 2011   # Stratford High Street
$ws.Range("B494").Value = 2009   # Stratford International
$ws.Range("B528").Value = 1987   # Tower Gateway
$ws.Range("B578").Value = 1987   # West India Quay
$ws.Range("B582").Value = 2005   # West Silvertown
$ws.Range("B587").Value = 1987   # Westferry
$ws.Range("B608").Value = 1849   # Woolwich Arsenal

# ------------------------------------------------------------------
# 3. Populate the DLR connections (station_from / station_to /
#    connectionYear) on the "Line_DLR" sheet.
# ------------------------------------------------------------------
$dlrWs = $wb.Worksheets.Item("Line_DLR")

$dlrWs.Range("H4").Value = "Bank"
$dlrWs.Range("I4").Value = "Shadwell"
$dlrWs.Range("J4").Value = 1994

$dlrWs.Range("H5").Value = "Tower Gateway"
$dlrWs.Range("I5").Value = "Shadwell"
$dlrWs.Range("J5").Value = 1994

$dlrWs.Range("I6").Value = "Limehouse"
$dlrWs.Range("J6").Value = 1994

$dlrWs.Range("I7").Value = "Westferry"
$dlrWs.Range("J7").Value = 1994

$dlrWs.Range("I8").Value = "Poplar"
$dlrWs.Range("J8").Value = 1994

$dlrWs.Range("I9").Value = "Blackwall"
$dlrWs.Range("J9").Value = 1994

$dlrWs.Range("I10").Value = "East India"
$dlrWs.Range("J10").Value = 1994

$dlrWs.Range("I11").Value = "Canning Town"
$dlrWs.Range("J11").Value = 1994

$dlrWs.Range("I12").Value = "Royal Victoria"
$dlrWs.Range("J12").Value = 1994

$dlrWs.Range("I13").Value = "Custom House"
$dlrWs.Range("J13").Value = 1994

$dlrWs.Range("I14").Value = "Prince Regent"
$dlrWs.Range("J14").Value = 1994

$dlrWs.Range("I15").Value = "Royal Albert"
$dlrWs.Range("J15").Value = 1994

$dlrWs.Range("I16").Value = "Beckton Park"
$dlrWs.Range("J16").Value = 1994

$dlrWs.Range("I17").Value = "Cyprus"
$dlrWs.Range("J17").Value = 1994

$dlrWs.Range("I18").Value = "Gallions Reach"
$dlrWs.Range("J18").Value = 1994

$dlrWs.Range("I19").Value = "Beckton"
$dlrWs.Range("J19").Value = 1994

$dlrWs.Range("H20").Value = "Canning Town"
$dlrWs.Range("I20").Value = "West Silvertown"
$dlrWs.Range("J20").Value = 1994

$dlrWs.Range("I21").Value = "Pontoon Dock"
$dlrWs.Range("J21").Value = 1994

$dlrWs.Range("I22").Value = "London City Airport"
$dlrWs.Range("J22").Value = 1994

$dlrWs.Range("I23").Value = "King George V"
$dlrWs.Range("J23").Value = 1994

$dlrWs.Range("I24").Value = "Woolwich Arsenal"
$dlrWs.Range("J24").Value = 1994

$dlrWs.Range("H25").Value = "Westferry"
$dlrWs.Range("I25").Value = "West India Quay"
$dlrWs.Range("J25").Value = 1994

$dlrWs.Range("I26").Value = "Canary Wharf"
$dlrWs.Range("J26").Value = 1994

$dlrWs.Range("I27").Value = "South Quay"
$dlrWs.Range("J27").Value = 1994

$dlrWs.Range("I28").Value = "Crossharbour"
$dlrWs.Range("J28").Value = 1994

$dlrWs.Range("I29").Value = "Mudchute"
$dlrWs.Range("J29").Value = 1994

$dlrWs.Range("I30").Value = "Island Gardens"
$dlrWs.Range("J30").Value = 1994

$dlrWs.Range("I31").Value = "Cutty Sark"
$dlrWs.Range("J31").Value = 1994

$dlrWs.Range("I32").Value = "Greenwich"
$dlrWs.Range("J32").Value = 1994

$dlrWs.Range("I33").Value = "Deptford Bridge"
$dlrWs.Range("J33").Value = 1994

$dlrWs.Range("I34").Value = "Elverson Road"
$dlrWs.Range("J34").Value = 1994

$dlrWs.Range("I35").Value = "Lewisham"
$dlrWs.Range("J35").Value = 1994

$dlrWs.Range("H36").Value = "Poplar"
$dlrWs.Range("I36").Value = "All Saints"
$dlrWs.Range("J36").Value = 1994

$dlrWs.Range("I37").Value = "Langdon Park"
$dlrWs.Range("J37").Value = 1994

$dlrWs.Range("I38").Value = "Devons Road"
$dlrWs.Range("J38").Value = 1994

$dlrWs.Range("I39").Value = "Bow Church"
$dlrWs.Range("J39").Value = 1994

$dlrWs.Range("I40").Value = "Pudding Mill Lane"
$dlrWs.Range("J40").Value = 1994

$dlrWs.Range("I41").Value = "Stratford"
$dlrWs.Range("J41").Value = 1994

$dlrWs.Range("H42").Value = "Canning Town"
$dlrWs.Range("I42").Value = "Star Lane"
$dlrWs.Range("J42").Value = 1994

$dlrWs.Range("I43").Value = "West Ham"
$dlrWs.Range("J43").Value = 1994

$dlrWs.Range("I44").Value = "Abbey Road"
$dlrWs.Range("J44").Value = 1994

$dlrWs.Range("I45").Value = "Stratford High Street"
$dlrWs.Range("J45").Value = 1994

$dlrWs.Range("I46").Value = "Stratford"
$dlrWs.Range("J46").Value = 1994

$dlrWs.Range("I47").Value = "Stratford International"
$dlrWs.Range("J47").Value = 1994

# ------------------------------------------------------------------
# 4. Column widths on the DLR sheet for the "station_from" (H) and
#    "station_to" (I) columns.
# ------------------------------------------------------------------
$dlrWs.Columns.Item(8).ColumnWidth = 14.5
$dlrWs.Columns.Item(9).ColumnWidth = 19.17

# ------------------------------------------------------------------
# 5. Restore / update the on-screen selections that were left behind
#    on the "Stations" and "Line_DLR" sheets.
# ------------------------------------------------------------------
$dlrWs.Activate() | Out-Null
$dlrWs.Range("L42").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("G314").Select() | Out-Null
